# Auto-generated Excel COM-interop script to apply scheduled-runner price updates
# across the Tonberry_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value2 = 107.5
$ws.Range("J33").Value2 = 100
$ws.Range("L33").Value2 = 100
$ws.Range("N33").Value2 = -558

$ws.Range("H40").Value2 = 2566.6667
$ws.Range("J40").Value2 = 1700
$ws.Range("L40").Value2 = 1700
$ws.Range("N40").Value2 = -2050

$ws.Range("H64").Value2 = 3048.889
$ws.Range("I64").Value2 = 2905.7144
$ws.Range("J64").Value2 = 3550
$ws.Range("K64").Value2 = 2905.7144
$ws.Range("L64").Value2 = 3550
$ws.Range("M64").Value2 = -2657.7144
$ws.Range("N64").Value2 = -4046

$ws.Range("H67").Value2 = 3048.889
$ws.Range("I67").Value2 = 2905.7144
$ws.Range("J67").Value2 = 3550
$ws.Range("K67").Value2 = 2905.7144
$ws.Range("L67").Value2 = 3550
$ws.Range("M67").Value2 = -2047.7144
$ws.Range("N67").Value2 = -5266

$ws.Range("H74").Value2 = 2996.6
$ws.Range("I74").Value2 = 2996.6
$ws.Range("K74").Value2 = 2996.6
$ws.Range("M74").Value2 = -2060.6

$ws.Range("H77").Value2 = 2996.6
$ws.Range("I77").Value2 = 2996.6
$ws.Range("K77").Value2 = 14983
$ws.Range("M77").Value2 = -10303

$ws.Range("H100").Value2 = 1883.2222
$ws.Range("I100").Value2 = 1868.625
$ws.Range("J100").Value2 = 2000
$ws.Range("K100").Value2 = 1868.625
$ws.Range("L100").Value2 = 2000
$ws.Range("M100").Value2 = -1327.625
$ws.Range("N100").Value2 = -3082

$ws.Range("H137").Value2 = 1383.1154
$ws.Range("I137").Value2 = 1235.9474
$ws.Range("K137").Value2 = 3707.8422
$ws.Range("M137").Value2 = -1157.8422

$ws.Range("H138").Value2 = 4030.1904
$ws.Range("I138").Value2 = 3816.2
$ws.Range("J138").Value2 = 4565.1665
$ws.Range("K138").Value2 = 11448.6
$ws.Range("L138").Value2 = 13695.4995
$ws.Range("M138").Value2 = -6308.599999999999
$ws.Range("N138").Value2 = -23975.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 2073.5
$ws.Range("I61").Value2 = 1149
$ws.Range("J61").Value2 = 4615.875
$ws.Range("K61").Value2 = 1149
$ws.Range("L61").Value2 = 4615.875
$ws.Range("M61").Value2 = -937
$ws.Range("N61").Value2 = -5039.875

$ws.Range("H132").Value2 = 1460.7838
$ws.Range("I132").Value2 = 982.4400000000001
$ws.Range("K132").Value2 = 2947.32
$ws.Range("M132").Value2 = -417.3200000000002

$ws.Range("H136").Value2 = 2073.5
$ws.Range("I136").Value2 = 1149
$ws.Range("J136").Value2 = 4615.875
$ws.Range("K136").Value2 = 3447
$ws.Range("L136").Value2 = 13847.625
$ws.Range("M136").Value2 = -897
$ws.Range("N136").Value2 = -18947.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 93691.55
$ws.Range("I86").Value2 = 2218.4546
$ws.Range("K86").Value2 = 2218.4546
$ws.Range("M86").Value2 = -1095.4546

$ws.Range("H89").Value2 = 93691.55
$ws.Range("I89").Value2 = 2218.4546
$ws.Range("K89").Value2 = 11092.273
$ws.Range("M89").Value2 = -5476.273000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value2 = 127.46667
$ws.Range("I7").Value2 = 69
$ws.Range("J7").Value2 = 244.4
$ws.Range("K7").Value2 = 69
$ws.Range("L7").Value2 = 244.4
$ws.Range("M7").Value2 = 44
$ws.Range("N7").Value2 = -470.4

$ws.Range("H16").Value2 = 802.5
$ws.Range("I16").Value2 = 760.4
$ws.Range("K16").Value2 = 760.4
$ws.Range("M16").Value2 = -473.4

$ws.Range("H31").Value2 = 1634.4445
$ws.Range("I31").Value2 = 999.2
$ws.Range("J31").Value2 = 2428.5
$ws.Range("K31").Value2 = 999.2
$ws.Range("L31").Value2 = 2428.5
$ws.Range("M31").Value2 = -704.2
$ws.Range("N31").Value2 = -3018.5

$ws.Range("H34").Value2 = 1634.4445
$ws.Range("I34").Value2 = 999.2
$ws.Range("J34").Value2 = 2428.5
$ws.Range("K34").Value2 = 999.2
$ws.Range("L34").Value2 = 2428.5
$ws.Range("M34").Value2 = -797.2
$ws.Range("N34").Value2 = -2832.5

$ws.Range("H44").Value2 = 3000
$ws.Range("J44").Value2 = 0
$ws.Range("L44").Value2 = 0
$ws.Range("N44").ClearContents()

$ws.Range("H50").Value2 = 21821.334
$ws.Range("J50").Value2 = 21821.334
$ws.Range("L50").Value2 = 21821.334
$ws.Range("N50").Value2 = -23071.334

$ws.Range("H58").Value2 = 1088539.6
$ws.Range("I58").Value2 = 1611465.9
$ws.Range("J58").Value2 = 2461.923
$ws.Range("K58").Value2 = 1611465.9
$ws.Range("L58").Value2 = 2461.923
$ws.Range("M58").Value2 = -1611262.9
$ws.Range("N58").Value2 = -2867.923

$ws.Range("H62").Value2 = 2673
$ws.Range("J62").Value2 = 2783.3333
$ws.Range("L62").Value2 = 2783.3333
$ws.Range("N62").Value2 = -4031.3333

$ws.Range("H65").Value2 = 2673
$ws.Range("J65").Value2 = 2783.3333
$ws.Range("L65").Value2 = 13916.6665
$ws.Range("N65").Value2 = -20156.6665

$ws.Range("H107").Value2 = 392.38095
$ws.Range("I107").Value2 = 297.5
$ws.Range("J107").Value2 = 696
$ws.Range("K107").Value2 = 297.5
$ws.Range("L107").Value2 = 696
$ws.Range("M107").Value2 = 1622.5
$ws.Range("N107").Value2 = -4536

$ws.Range("H113").Value2 = 802.5
$ws.Range("I113").Value2 = 760.4
$ws.Range("K113").Value2 = 760.4
$ws.Range("M113").Value2 = 1409.6

$ws.Range("H132").Value2 = 1648.9491
$ws.Range("I132").Value2 = 1118.925
$ws.Range("K132").Value2 = 3356.775
$ws.Range("M132").Value2 = -826.7749999999996

$ws.Range("H136").Value2 = 1088539.6
$ws.Range("I136").Value2 = 1611465.9
$ws.Range("J136").Value2 = 2461.923
$ws.Range("K136").Value2 = 4834397.699999999
$ws.Range("L136").Value2 = 7385.768999999999
$ws.Range("M136").Value2 = -4831847.699999999
$ws.Range("N136").Value2 = -12485.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value2 = 991.75
$ws.Range("J122").Value2 = 1225.8572
$ws.Range("L122").Value2 = 11032.7148
$ws.Range("N122").Value2 = -15932.7148

$ws.Range("H131").Value2 = 789.95
$ws.Range("I131").Value2 = 377.22223
$ws.Range("J131").Value2 = 830.7692
$ws.Range("K131").Value2 = 1131.66669
$ws.Range("L131").Value2 = 2492.3076
$ws.Range("M131").Value2 = 3908.33331
$ws.Range("N131").Value2 = -12572.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 1204420.2
$ws.Range("I132").Value2 = 1925137.2
$ws.Range("J132").Value2 = 3225.0833
$ws.Range("K132").Value2 = 5775411.6
$ws.Range("L132").Value2 = 9675.249899999999
$ws.Range("M132").Value2 = -5772881.6
$ws.Range("N132").Value2 = -14735.2499

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 3246.3076
$ws.Range("I7").Value2 = 2137
$ws.Range("K7").Value2 = 2137
$ws.Range("M7").Value2 = -2025

$ws.Range("H46").Value2 = 2487.4
$ws.Range("I46").Value2 = 1425
$ws.Range("K46").Value2 = 1425
$ws.Range("M46").Value2 = -1237

$ws.Range("H82").Value2 = 2130.5334
$ws.Range("I82").Value2 = 1570
$ws.Range("J82").Value2 = 2621
$ws.Range("K82").Value2 = 1570
$ws.Range("L82").Value2 = 2621
$ws.Range("M82").Value2 = -1209
$ws.Range("N82").Value2 = -3343

$ws.Range("H85").Value2 = 2130.5334
$ws.Range("I85").Value2 = 1570
$ws.Range("J85").Value2 = 2621
$ws.Range("K85").Value2 = 1570
$ws.Range("L85").Value2 = 2621
$ws.Range("M85").Value2 = -322
$ws.Range("N85").Value2 = -5117

$ws.Range("H122").Value2 = 9661.6
$ws.Range("J122").Value2 = 14966.667
$ws.Range("L122").Value2 = 44900.001
$ws.Range("N122").Value2 = -49800.001

$ws.Range("H126").Value2 = 3246.3076
$ws.Range("I126").Value2 = 2137
$ws.Range("K126").Value2 = 6411
$ws.Range("M126").Value2 = -3941

$ws.Range("H132").Value2 = 1657.8
$ws.Range("J132").Value2 = 1980.1364
$ws.Range("L132").Value2 = 5940.4092
$ws.Range("N132").Value2 = -11000.4092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 1678.7273
$ws.Range("I81").Value2 = 1546.6
$ws.Range("K81").Value2 = 3093.2
$ws.Range("M81").Value2 = -2032.2

$ws.Range("H84").Value2 = 1678.7273
$ws.Range("I84").Value2 = 1546.6
$ws.Range("K84").Value2 = 15466
$ws.Range("M84").Value2 = -10162

$ws.Range("H96").Value2 = 3666.6667
$ws.Range("I96").Value2 = 3000
$ws.Range("K96").Value2 = 3000
$ws.Range("M96").Value2 = -1627
